$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed shared-string order: "pk" needs to land before "g" in sharedStrings.xml
$ws.Range("K39").Value = "pk"
$ws.Range("K36").Value = "g"

# Row 36 - Simcoe hop addition
$ws.Range("I36").Value = "Simcoe"
$ws.Range("J36").Value = 100
$ws.Range("K36").Value = "g"

# Row 37 - Cascade hop addition
$ws.Range("I37").Value = "Cascade"
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = "g"

# Row 38 - Chinook hop addition
$ws.Range("I38").Value = "Chinook"
$ws.Range("J38").Value = 100
$ws.Range("K38").Value = "g"

# Row 39 - Gjaer (yeast) addition
$ws.Range("I39").Value = "Gjær"
$ws.Range("J39").Value = 2
$ws.Range("K39").Value = "pk"

# Update the view: scroll and select cell Q38
$ws.Range("Q38").Select()
